$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 151, shifting existing rows 151:161 down to 152:162
$ws.Rows("151:151").Insert()

# Populate the new row 151 with the new weekly price entry
$ws.Range("A151").Value = 1
$ws.Range("B151").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C151").Value = "Arica y Parinacota"
$ws.Range("D151").Value = 44491
$ws.Range("D151").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E151").Value = 15
$ws.Range("F151").Value = "Fruta"
$ws.Range("G151").Value = 100102
$ws.Range("H151").Value = "Cítricos"
$ws.Range("I151").Value = 100102003
$ws.Range("J151").Value = "Limón"
$ws.Range("K151").Value = "Tahití"
$ws.Range("L151").Value = "Primera"
$ws.Range("M151").Value = 200
$ws.Range("N151").Value = 45000
$ws.Range("O151").Value = 46000
$ws.Range("P151").Value = 45500
$ws.Range("Q151").Value = "`$/caja 24 kilos"
$ws.Range("R151").Value = "Perú"
$ws.Range("S151").Value = 1896
$ws.Range("T151").Value = 24
